$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.412.36"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "2.776.02"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "352.52"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").Value = "107.94"
$ws.Range("D7").Value = "0.549"
$ws.Range("E7").Value = "  -2.75%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").Value = "39.66"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("E11").Value = "  +3.09%  "
$ws.Range("D12").Value = "0.0836"
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("D13").Value = "20.04"
$ws.Range("E13").Value = "  +3.11%  "
$ws.Range("D14").Value = "7.62"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "3.211.49"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "2.761.69"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").Value = "0.923"
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").Value = "51.397.58"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "7.64"
$ws.Range("E19").Value = "  +2.87%  "
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").Value = "13.12"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").Value = "69.87"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "265.63"
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("D25").Value = "2.70"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "26.04"
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("D28").Value = "0.163"
$ws.Range("E28").Value = "  +12.36%  "
$ws.Range("D29").Value = "10.22"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").Value = "36.20"
$ws.Range("E31").Value = "  +6.46%  "
$ws.Range("D32").Value = "6.25"
$ws.Range("E32").Value = "  +9.61%  "
$ws.Range("D33").Value = "51.93"
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("D34").Value = "0.0454"
$ws.Range("E34").Value = "  -2.48%  "
$ws.Range("E35").Value = "  +5.88%  "
$ws.Range("E36").Value = "  -2.30%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "18.45"
$ws.Range("E38").Value = "  +1.86%  "
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "22.10"
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "120.04"
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("E45").Value = "  -2.17%  "
$ws.Range("D46").Value = "2.106.39"
$ws.Range("E46").Value = "  +2.01%  "
$ws.Range("D47").Value = "3.26"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("E48").Value = "  +6.23%  "
$ws.Range("B49").Value = "SEI"
$ws.Range("C49").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D49").Value = "0.906"
$ws.Range("E49").Value = "  -2.91%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "5.41"
$ws.Range("E50").Value = "  -4.88%  "
$ws.Range("E51").Value = "  +7.37%  "

Write-Output "done"
